# This script reorders the "Requisitos" course list so that it is sorted
# alphabetically by course code (LOB1003, LOB1004, LOB1006, ... LOQ4251).
# Each course line lives in its own run ("<code> -  <name>  (Requisito)" + a
# line break) inside one ListBullet paragraph. Because every one of the 19
# lines needs to move, we cannot safely Find/Replace old-text -> new-text in a
# single pass: the 19 "new" strings are exactly the same 19 strings as the
# "old" ones (just reordered), so a naive pass could match text that a later
# step already placed. We avoid this by going through unique placeholder
# tokens first (phase 1), then resolving the placeholders to the final,
# sorted text (phase 2).

$d = $word.ActiveDocument

$oldOrder = @(
    "LOB1045 -  Leitura e Produção de Textos Acadêmicos  (Requisito)",
    "LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito)",
    "LOQ4251 -  Fundamentos de Química  (Requisito)",
    "LOB1006 -  Cálculo IV  (Requisito)",
    "LOB1037 -  Àlgebra Linear  (Requisito)",
    "LOB1053 -  Física III  (Requisito)",
    "LOB1003 -  Cálculo I  (Requisito)",
    "LOB1009 -  Leitura e Interpretação de Desenho Técnico  (Requisito)",
    "LOB1012 -  Estatística  (Requisito)",
    "LOB1018 -  Física I  (Requisito)",
    "LOB1024 -  Mecânica  (Requisito)",
    "LOB1036 -  Geometria Analítica  (Requisito)",
    "LOB1038 -  Física Experimental I  (Requisito)",
    "LOB1039 -  Física Experimental III  (Requisito)",
    "LOB1041 -  Física Experimental II  (Requisito)",
    "LOB1052 -  Cálculo III  (Requisito)",
    "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)",
    "LOB1004 -  Cálculo II  (Requisito)",
    "LOB1019 -  Física II  (Requisito)"
)

$newOrder = @(
    "LOB1003 -  Cálculo I  (Requisito)",
    "LOB1004 -  Cálculo II  (Requisito)",
    "LOB1006 -  Cálculo IV  (Requisito)",
    "LOB1009 -  Leitura e Interpretação de Desenho Técnico  (Requisito)",
    "LOB1012 -  Estatística  (Requisito)",
    "LOB1018 -  Física I  (Requisito)",
    "LOB1019 -  Física II  (Requisito)",
    "LOB1024 -  Mecânica  (Requisito)",
    "LOB1036 -  Geometria Analítica  (Requisito)",
    "LOB1037 -  Àlgebra Linear  (Requisito)",
    "LOB1038 -  Física Experimental I  (Requisito)",
    "LOB1039 -  Física Experimental III  (Requisito)",
    "LOB1041 -  Física Experimental II  (Requisito)",
    "LOB1045 -  Leitura e Produção de Textos Acadêmicos  (Requisito)",
    "LOB1052 -  Cálculo III  (Requisito)",
    "LOB1053 -  Física III  (Requisito)",
    "LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito)",
    "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)",
    "LOQ4251 -  Fundamentos de Química  (Requisito)"
)

# Phase 1: old text -> unique placeholder token (order of the tokens does not matter,
# only that each is unique and cannot collide with any real course text).
for ($i = 0; $i -lt $oldOrder.Length; $i++) {
    $placeholder = "@@REQ_SLOT_{0:D2}@@" -f $i
    $d.Content.Find.Execute($oldOrder[$i], $true, $false, $false, $false, $false, $true, 1, $false, $placeholder, 2) | Out-Null
}

# Phase 2: placeholder token -> final text, in the sorted order.
for ($i = 0; $i -lt $newOrder.Length; $i++) {
    $placeholder = "@@REQ_SLOT_{0:D2}@@" -f $i
    $d.Content.Find.Execute($placeholder, $true, $false, $false, $false, $false, $true, 1, $false, $newOrder[$i], 2) | Out-Null
}
